# Presentacion semana III 05 de Julio -> "Actualización de dpcumentos, Semana 4"
#
# 1) The cached "last updated" date shown by the datetimeFigureOut field
#    placeholder (slide master + every slide layout) moves from
#    9/4/2013 to 9/5/2013.
# 2) The small status table on the slide ("Ultima Revisión" / "Documentos"
#    / "Aplicación") has its two numeric cells bumped by one (13->14,
#    14->15).

$p = $ppt.ActivePresentation

$oldDate = "9/4/2013"
$newDate = "9/5/2013"

# --- 1) Slide master date placeholder -------------------------------------
for ($i = 1; $i -le $p.SlideMaster.Shapes.Count; $i++) {
    $sh = $p.SlideMaster.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 1b) Every custom (slide) layout's date placeholder --------------------
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $sh = $lay.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# --- 2) Status table numbers -------------------------------------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $tbl = $sh.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cellShape = $tbl.Cell($r, $c).Shape
                if ($cellShape.TextFrame.HasText) {
                    $t = $cellShape.TextFrame.TextRange.Text
                    if ($t -eq "13") {
                        $cellShape.TextFrame.TextRange.Text = "14"
                    } elseif ($t -eq "14") {
                        $cellShape.TextFrame.TextRange.Text = "15"
                    }
                }
            }
        }
    }
}
